# Updated cryptos list on Fri Sep 29 21:32:49 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.917.27'
$ws.Range("E2").Value = '  -0.29%  '

$ws.Range("D3").Value = '1.669.09'
$ws.Range("E3").Value = '  +1.19%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").Value = '''215.82'
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("D8").Value = '''0.0621'
$ws.Range("E8").Value = '  +1.29%  '

$ws.Range("E9").Value = '  +0.34%  '

$ws.Range("D10").Value = '''20.28'
$ws.Range("E10").Value = '  +3.14%  '

$ws.Range("D11").Value = '''0.0895'
$ws.Range("E11").Value = '  +3.44%  '

$ws.Range("D12").Value = '1.903.15'
$ws.Range("E12").Value = '  +1.05%  '

$ws.Range("D13").Value = '1.677.29'
$ws.Range("E13").Value = '  +1.66%  '

$ws.Range("E14").Value = '  +0.44%  '

$ws.Range("D15").Value = '''0.526'
$ws.Range("E15").Value = '  +1.73%  '

$ws.Range("D16").Value = '''66.19'
$ws.Range("E16").Value = '  +1.68%  '

$ws.Range("D17").Value = '26.942.54'
$ws.Range("E17").Value = '  -0.14%  '

$ws.Range("D18").Value = '''234.52'
$ws.Range("E18").Value = '  -1.36%  '

$ws.Range("D19").Value = '''7.96'
$ws.Range("E19").Value = '  +1.85%  '

$ws.Range("E20").Value = '  +0.58%  '

$ws.Range("E21").Value = '  -0.12%  '

$ws.Range("D22").Value = '''4.40'
$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("E23").Value = '  -2.18%  '

$ws.Range("E24").Value = '  -0.82%  '

$ws.Range("D25").Value = '''146.09'
$ws.Range("E25").Value = '  +0.25%  '

$ws.Range("D26").Value = '''7.15'
$ws.Range("E26").Value = '  +0.70%  '

$ws.Range("E27").Value = '  +1.23%  '

$ws.Range("D28").Value = '''15.93'
$ws.Range("E28").Value = '  +0.76%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("E30").Value = '  +0.17%  '

$ws.Range("E31").Value = '  +0.26%  '

$ws.Range("D32").Value = '''3.36'
$ws.Range("E32").Value = '  +2.10%  '

$ws.Range("D33").Value = '1.452.19'
$ws.Range("E33").Value = '  -3.84%  '

$ws.Range("E34").Value = '  +2.50%  '

$ws.Range("E35").Value = '  +4.60%  '

$ws.Range("E36").Value = '  -0.50%  '

$ws.Range("D37").Value = '''0.583'
$ws.Range("E37").Value = '  +1.64%  '

$ws.Range("E38").Value = '  +2.28%  '

$ws.Range("E39").Value = '  +0.74%  '

$ws.Range("E40").Value = '  -3.46%  '

$ws.Range("E41").Value = '  -0.17%  '

$ws.Range("E42").Value = '  +1.42%  '

$ws.Range("D43").Value = '''66.20'
$ws.Range("E43").Value = '  +0.53%  '

$ws.Range("D44").Value = '''0.974'
$ws.Range("E44").Value = '  +6.11%  '

$ws.Range("D45").Value = '1.811.07'
$ws.Range("E45").Value = '  +1.09%  '

$ws.Range("D46").Value = '''0.784'
$ws.Range("E46").Value = '  +1.25%  '

$ws.Range("D47").Value = '''90.64'
$ws.Range("E47").Value = '  +1.42%  '

$ws.Range("E48").Value = '  +1.55%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0104'
$ws.Range("E49").Value = '  -0.38%  '

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '''0.102'
$ws.Range("E50").Value = '  +4.87%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.0505'
$ws.Range("E51").Value = '  -0.30%  '
